# Refresh the "cryptos" price table: update Price (D), Volume(1h) (E), and
# Hora (G, hour-of-day) columns for rows 2-51 per the latest GitHub Actions
# scrape. All values are stored as text (matching the sheet's existing
# inline/shared-string convention), so a leading apostrophe is used to force
# Excel to keep them as text instead of auto-converting to numbers/percentages
# (which would also mangle small decimals like 0.00006243 into scientific
# notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.46"
$ws.Range("E2").Value = "'-1.22%"
$ws.Range("G2").Value = "'15"
$ws.Range("D3").Value = "'31.45"
$ws.Range("E3").Value = "'-1.75%"
$ws.Range("G3").Value = "'15"
$ws.Range("D4").Value = "'5.145"
$ws.Range("E4").Value = "'-2.97%"
$ws.Range("G4").Value = "'15"
$ws.Range("E5").Value = "'-1.64%"
$ws.Range("G5").Value = "'15"
$ws.Range("D6").Value = "'1.835"
$ws.Range("E6").Value = "'25.01%"
$ws.Range("G6").Value = "'15"
$ws.Range("D7").Value = "'7.784"
$ws.Range("E7").Value = "'-0.75%"
$ws.Range("G7").Value = "'15"
$ws.Range("D8").Value = "'3.751"
$ws.Range("E8").Value = "'-0.54%"
$ws.Range("G8").Value = "'15"
$ws.Range("D9").Value = "'0.9257"
$ws.Range("E9").Value = "'0.80%"
$ws.Range("G9").Value = "'15"
$ws.Range("E10").Value = "'-0.49%"
$ws.Range("G10").Value = "'15"
$ws.Range("D11").Value = "'0.07089"
$ws.Range("E11").Value = "'-8.10%"
$ws.Range("G11").Value = "'15"
$ws.Range("D12").Value = "'0.08045"
$ws.Range("E12").Value = "'-0.45%"
$ws.Range("G12").Value = "'15"
$ws.Range("D13").Value = "'0.03033"
$ws.Range("E13").Value = "'0.35%"
$ws.Range("G13").Value = "'15"
$ws.Range("D14").Value = "'0.09922"
$ws.Range("E14").Value = "'0.40%"
$ws.Range("G14").Value = "'15"
$ws.Range("D15").Value = "'0.001498"
$ws.Range("E15").Value = "'0.21%"
$ws.Range("G15").Value = "'15"
$ws.Range("D16").Value = "'0.006154"
$ws.Range("E16").Value = "'-1.08%"
$ws.Range("G16").Value = "'15"
$ws.Range("E17").Value = "'-0.70%"
$ws.Range("G17").Value = "'15"
$ws.Range("D18").Value = "'2.222"
$ws.Range("E18").Value = "'-0.29%"
$ws.Range("G18").Value = "'15"
$ws.Range("G19").Value = "'15"
$ws.Range("D20").Value = "'0.1330"
$ws.Range("E20").Value = "'-1.00%"
$ws.Range("G20").Value = "'15"
$ws.Range("D21").Value = "'4.560"
$ws.Range("E21").Value = "'1.19%"
$ws.Range("G21").Value = "'15"
$ws.Range("D22").Value = "'0.04649"
$ws.Range("E22").Value = "'1.76%"
$ws.Range("G22").Value = "'15"
$ws.Range("E23").Value = "'-2.69%"
$ws.Range("G23").Value = "'15"
$ws.Range("D24").Value = "'0.001214"
$ws.Range("E24").Value = "'-0.31%"
$ws.Range("G24").Value = "'15"
$ws.Range("D25").Value = "'0.004746"
$ws.Range("E25").Value = "'7.39%"
$ws.Range("G25").Value = "'15"
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("E26").Value = "'-7.41%"
$ws.Range("G26").Value = "'15"
$ws.Range("E27").Value = "'7.55%"
$ws.Range("G27").Value = "'15"
$ws.Range("G28").Value = "'15"
$ws.Range("G29").Value = "'15"
$ws.Range("G30").Value = "'15"
$ws.Range("G31").Value = "'15"
$ws.Range("G32").Value = "'15"
$ws.Range("G33").Value = "'15"
$ws.Range("G34").Value = "'15"
$ws.Range("G35").Value = "'15"
$ws.Range("G36").Value = "'15"
$ws.Range("G37").Value = "'15"
$ws.Range("G38").Value = "'15"
$ws.Range("D39").Value = "'0.01721"
$ws.Range("E39").Value = "'-0.44%"
$ws.Range("G39").Value = "'15"
$ws.Range("D40").Value = "'0.04496"
$ws.Range("E40").Value = "'-0.77%"
$ws.Range("G40").Value = "'15"
$ws.Range("D41").Value = "'0.007107"
$ws.Range("E41").Value = "'-1.41%"
$ws.Range("G41").Value = "'15"
$ws.Range("E42").Value = "'-0.38%"
$ws.Range("G42").Value = "'15"
$ws.Range("D43").Value = "'0.002176"
$ws.Range("E43").Value = "'-2.09%"
$ws.Range("G43").Value = "'15"
$ws.Range("D44").Value = "'0.01085"
$ws.Range("E44").Value = "'-20.27%"
$ws.Range("G44").Value = "'15"
$ws.Range("D45").Value = "'0.00006243"
$ws.Range("E45").Value = "'1.55%"
$ws.Range("G45").Value = "'15"
$ws.Range("E46").Value = "'-21.39%"
$ws.Range("G46").Value = "'15"
$ws.Range("G47").Value = "'15"
$ws.Range("G48").Value = "'15"
$ws.Range("G49").Value = "'15"
$ws.Range("G50").Value = "'15"
$ws.Range("G51").Value = "'15"

# Reset style on touched data cells (rows 2-51) back to the default/unstyled
# appearance, since assigning a leading apostrophe (quote-prefix) to force
# text storage also tags the cell with a quotePrefix style.
$ws.Range("D2:E51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
